$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 138 (pushes the old rows 138:217 down to 139:218,
# extending the used range to A1:R218) and fill it with the new weekly record.
$ws.Rows.Item(138).Insert()

$ws.Range("A138").Value = 9
$ws.Range("B138").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C138").Value = 'Metropolitana'
$ws.Range("D138").Value = 44596
$ws.Range("E138").Value = 13
$ws.Range("F138").Value = 300000001
$ws.Range("G138").Value = 'Rabanito'
$ws.Range("H138").Value = 'Sin especificar'
$ws.Range("I138").Value = 'Primera'
$ws.Range("J138").Value = 14000
$ws.Range("K138").Value = 3500
$ws.Range("L138").Value = 4000
$ws.Range("M138").Value = 3786
$ws.Range("N138").Value = '$/cien unidades (volumen en unidades)'
$ws.Range("O138").Value = 'Provincia de Chacabuco'
$ws.Range("P138").Value = 38
$ws.Range("Q138").Value = 100
$ws.Range("R138").Value = 'Hortaliza'
